$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_DM")

# Insert two new blank columns before column X (24), shifting codelist..change_history right by two
$ws.Range("X1:Y1").EntireColumn.Insert()

# Set the new header values
$ws.Range("X1").Value = "derived_variable"
$ws.Range("Y1").Value = "derivation_description"

# Update the view: scroll so column R is the top-left visible column, and select Z14
$ws.Application.ActiveWindow.ScrollColumn = 18
$ws.Range("Z14").Select()

$wb.Save()
